# The worksheet that actually carries this test data is "Add Panels and
# Devices" (the first tab) - make sure we operate on that one explicitly
# rather than relying on whatever happens to be the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Add Panels and Devices")

# F8 used to hold the numeric value 0.393; it is now entered as text that
# looks like a number (typed with a leading apostrophe in Excel), which is
# why it ends up stored as a shared string together with a "quote prefix"
# cell style.
$ws.Range("F8").Value2 = "'0.329"

# J8 and K8 keep their numeric type but get new values.
$ws.Range("J8").Value2 = 0.405
$ws.Range("K8").Value2 = 0.329

# Reflect the new selection/active cell on this sheet.
$ws.Activate()
$ws.Range("F8").Select()
